$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 430.125
$ws.Range("J12").Value = 737.5
$ws.Range("L12").Value = 737.5
$ws.Range("N12").Value = -1077.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 799.5
$ws.Range("I32").Value = 799
$ws.Range("K32").Value = 799
$ws.Range("M32").Value = -473

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 128.15384
$ws.Range("I33").Value = 90.666664
$ws.Range("K33").Value = 90.666664
$ws.Range("M33").Value = 138.333336

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2831.6667
$ws.Range("J51").Value = 2898
$ws.Range("L51").Value = 2898
$ws.Range("N51").Value = -3866

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 251.25
$ws.Range("J92").Value = 300
$ws.Range("L92").Value = 300
$ws.Range("N92").Value = -2796

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4999.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8339251
$ws.Range("I32").Value = 2251.75
$ws.Range("K32").Value = 2251.75
$ws.Range("M32").Value = -1964.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1472.875
$ws.Range("I61").Value = 1501
$ws.Range("J61").Value = 1388.5
$ws.Range("K61").Value = 1501
$ws.Range("L61").Value = 1388.5
$ws.Range("M61").Value = -1289
$ws.Range("N61").Value = -1812.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2679.25
$ws.Range("I97").Value = 2995
$ws.Range("K97").Value = 2995
$ws.Range("M97").Value = -2499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2513.5386
$ws.Range("I110").Value = 1964.4286
$ws.Range("K110").Value = 1964.4286
$ws.Range("M110").Value = 80.57140000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2879.7144
$ws.Range("J122").Value = 2582
$ws.Range("L122").Value = 7746
$ws.Range("N122").Value = -12646

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2562.25
$ws.Range("I132").Value = 2714.75
$ws.Range("J132").Value = 1799.75
$ws.Range("K132").Value = 8144.25
$ws.Range("L132").Value = 5399.25
$ws.Range("M132").Value = -5614.25
$ws.Range("N132").Value = -10459.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1472.875
$ws.Range("I136").Value = 1501
$ws.Range("J136").Value = 1388.5
$ws.Range("K136").Value = 4503
$ws.Range("L136").Value = 4165.5
$ws.Range("M136").Value = -1953
$ws.Range("N136").Value = -9265.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 573.375
$ws.Range("I94").Value = 497.83334
$ws.Range("K94").Value = 497.83334
$ws.Range("M94").Value = -46.83334000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1951.25
$ws.Range("I99").Value = 1955.909
$ws.Range("K99").Value = 1955.909
$ws.Range("M99").Value = -457.9090000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5805.121
$ws.Range("I31").Value = 3532.5557
$ws.Range("J31").Value = 6657.3335
$ws.Range("K31").Value = 3532.5557
$ws.Range("L31").Value = 6657.3335
$ws.Range("M31").Value = -3237.5557
$ws.Range("N31").Value = -7247.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5805.121
$ws.Range("I34").Value = 3532.5557
$ws.Range("J34").Value = 6657.3335
$ws.Range("K34").Value = 3532.5557
$ws.Range("L34").Value = 6657.3335
$ws.Range("M34").Value = -3330.5557
$ws.Range("N34").Value = -7061.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1923.3334
$ws.Range("I105").Value = 1507.25
$ws.Range("J105").Value = 2755.5
$ws.Range("K105").Value = 1507.25
$ws.Range("L105").Value = 2755.5
$ws.Range("M105").Value = 239.75
$ws.Range("N105").Value = -6249.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1156.8
$ws.Range("I107").Value = 937.9
$ws.Range("J107").Value = 1594.6
$ws.Range("K107").Value = 937.9
$ws.Range("L107").Value = 1594.6
$ws.Range("M107").Value = 982.1
$ws.Range("N107").Value = -5434.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2650
$ws.Range("J5").Value = 2675
$ws.Range("L5").Value = 8025
$ws.Range("N5").Value = -8249

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 823.5
$ws.Range("I122").Value = 873.8333
$ws.Range("K122").Value = 7864.4997
$ws.Range("M122").Value = -5414.4997

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 2650
$ws.Range("J135").Value = 2675
$ws.Range("L135").Value = 24075
$ws.Range("N135").Value = -29145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 5519.125
$ws.Range("I55").Value = 5859
$ws.Range("J55").Value = 4499.5
$ws.Range("K55").Value = 5859
$ws.Range("L55").Value = 4499.5
$ws.Range("M55").Value = -5532
$ws.Range("N55").Value = -5153.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2959.8
$ws.Range("I70").Value = 2959.8
$ws.Range("K70").Value = 2959.8
$ws.Range("M70").Value = -2689.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 2959.8
$ws.Range("I73").Value = 2959.8
$ws.Range("K73").Value = 2959.8
$ws.Range("M73").Value = -2023.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 915.25
$ws.Range("J97").Value = 860.5
$ws.Range("L97").Value = 860.5
$ws.Range("N97").Value = -1852.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1534.9231
$ws.Range("I102").Value = 1534.9231
$ws.Range("K102").Value = 1534.9231
$ws.Range("M102").Value = 87.07690000000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 62804.766
$ws.Range("I132").Value = 93548.37
$ws.Range("J132").Value = 6441.5
$ws.Range("K132").Value = 280645.11
$ws.Range("L132").Value = 19324.5
$ws.Range("M132").Value = -278115.11
$ws.Range("N132").Value = -24384.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4438.125
$ws.Range("J46").Value = 4706.9165
$ws.Range("L46").Value = 4706.9165
$ws.Range("N46").Value = -5082.9165

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 3985
$ws.Range("I55").Value = 3985
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 3985
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -3812
$ws.Range("N55").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 7000
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 7000
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 7000
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -7404

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2911.3333
$ws.Range("I82").Value = 1115.1666
$ws.Range("K82").Value = 1115.1666
$ws.Range("M82").Value = -754.1666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2911.3333
$ws.Range("I85").Value = 1115.1666
$ws.Range("K85").Value = 1115.1666
$ws.Range("M85").Value = 132.8334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1347.7
$ws.Range("I93").Value = 1336.5333
$ws.Range("K93").Value = 1336.5333
$ws.Range("M93").Value = -88.53330000000005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 7000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 7000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 7000
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -11340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 69500.5
$ws.Range("J127").Value = 69500.5
$ws.Range("L127").Value = 69500.5
$ws.Range("N127").Value = -79420.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2612.8572
$ws.Range("I136").Value = 2465
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 7395
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -4845
$ws.Range("N136").Value = -15600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 44999
$ws.Range("J47").Value = 44999
$ws.Range("L47").Value = 44999
$ws.Range("N47").Value = -46143

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8141.0835
$ws.Range("J62").Value = 8862.777
$ws.Range("L62").Value = 8862.777
$ws.Range("N62").Value = -10110.777

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 8141.0835
$ws.Range("J65").Value = 8862.777
$ws.Range("L65").Value = 44313.885
$ws.Range("N65").Value = -50553.885

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1113.4166
$ws.Range("J96").Value = 867.6
$ws.Range("L96").Value = 867.6
$ws.Range("N96").Value = -3613.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1758.5
$ws.Range("I100").Value = 1809.0667
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 3618.1334
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -3077.1334
$ws.Range("N100").Value = -3082
